$wb = $excel.ActiveWorkbook

# ---- PIR sheet: rows 177-189 ----
$wsPIR = $wb.Worksheets.Item("PIR")
$wsPIR.Range('A177').NumberFormat = '@'
$wsPIR.Range('A177').Value = '2026-01-28'
$wsPIR.Range('B177').Value = '18:16:07'
$wsPIR.Range('C177').Value = '18:00'
$wsPIR.Range('D177').Value = 'Bathroom'
$wsPIR.Range('E177').Value = 'No Motion'
$wsPIR.Range('F177').Value = 'Inactive'

$wsPIR.Range('A178').NumberFormat = '@'
$wsPIR.Range('A178').Value = '2026-01-28'
$wsPIR.Range('B178').Value = '18:16:10'
$wsPIR.Range('C178').Value = '18:00'
$wsPIR.Range('D178').Value = 'Bathroom'
$wsPIR.Range('E178').Value = 'No Motion'
$wsPIR.Range('F178').Value = 'Inactive'

$wsPIR.Range('A179').NumberFormat = '@'
$wsPIR.Range('A179').Value = '2026-01-28'
$wsPIR.Range('B179').Value = '18:16:13'
$wsPIR.Range('C179').Value = '18:00'
$wsPIR.Range('D179').Value = 'Bathroom'
$wsPIR.Range('E179').Value = 'No Motion'
$wsPIR.Range('F179').Value = 'Inactive'

$wsPIR.Range('A180').NumberFormat = '@'
$wsPIR.Range('A180').Value = '2026-01-28'
$wsPIR.Range('B180').Value = '18:16:18'
$wsPIR.Range('C180').Value = '18:00'
$wsPIR.Range('D180').Value = 'Bathroom'
$wsPIR.Range('E180').Value = 'No Motion'
$wsPIR.Range('F180').Value = 'Inactive'

$wsPIR.Range('A181').NumberFormat = '@'
$wsPIR.Range('A181').Value = '2026-01-28'
$wsPIR.Range('B181').Value = '18:16:22'
$wsPIR.Range('C181').Value = '18:00'
$wsPIR.Range('D181').Value = 'Bathroom'
$wsPIR.Range('E181').Value = 'No Motion'
$wsPIR.Range('F181').Value = 'Inactive'

$wsPIR.Range('A182').NumberFormat = '@'
$wsPIR.Range('A182').Value = '2026-01-28'
$wsPIR.Range('B182').Value = '18:16:26'
$wsPIR.Range('C182').Value = '18:00'
$wsPIR.Range('D182').Value = 'Bathroom'
$wsPIR.Range('E182').Value = 'No Motion'
$wsPIR.Range('F182').Value = 'Inactive'

$wsPIR.Range('A183').NumberFormat = '@'
$wsPIR.Range('A183').Value = '2026-01-28'
$wsPIR.Range('B183').Value = '18:16:31'
$wsPIR.Range('C183').Value = '18:00'
$wsPIR.Range('D183').Value = 'Bathroom'
$wsPIR.Range('E183').Value = 'No Motion'
$wsPIR.Range('F183').Value = 'Inactive'

$wsPIR.Range('A184').NumberFormat = '@'
$wsPIR.Range('A184').Value = '2026-01-28'
$wsPIR.Range('B184').Value = '18:16:36'
$wsPIR.Range('C184').Value = '18:00'
$wsPIR.Range('D184').Value = 'Bathroom'
$wsPIR.Range('E184').Value = 'No Motion'
$wsPIR.Range('F184').Value = 'Inactive'

$wsPIR.Range('A185').NumberFormat = '@'
$wsPIR.Range('A185').Value = '2026-01-28'
$wsPIR.Range('B185').Value = '18:16:41'
$wsPIR.Range('C185').Value = '18:00'
$wsPIR.Range('D185').Value = 'Bathroom'
$wsPIR.Range('E185').Value = 'No Motion'
$wsPIR.Range('F185').Value = 'Inactive'

$wsPIR.Range('A186').NumberFormat = '@'
$wsPIR.Range('A186').Value = '2026-01-28'
$wsPIR.Range('B186').Value = '18:16:47'
$wsPIR.Range('C186').Value = '18:00'
$wsPIR.Range('D186').Value = 'Bathroom'
$wsPIR.Range('E186').Value = 'No Motion'
$wsPIR.Range('F186').Value = 'Inactive'

$wsPIR.Range('A187').NumberFormat = '@'
$wsPIR.Range('A187').Value = '2026-01-28'
$wsPIR.Range('B187').Value = '18:16:52'
$wsPIR.Range('C187').Value = '18:00'
$wsPIR.Range('D187').Value = 'Bathroom'
$wsPIR.Range('E187').Value = 'No Motion'
$wsPIR.Range('F187').Value = 'Inactive'

$wsPIR.Range('A188').NumberFormat = '@'
$wsPIR.Range('A188').Value = '2026-01-28'
$wsPIR.Range('B188').Value = '18:16:59'
$wsPIR.Range('C188').Value = '18:00'
$wsPIR.Range('D188').Value = 'Bathroom'
$wsPIR.Range('E188').Value = 'No Motion'
$wsPIR.Range('F188').Value = 'Inactive'

$wsPIR.Range('A189').NumberFormat = '@'
$wsPIR.Range('A189').Value = '2026-01-28'
$wsPIR.Range('B189').Value = '18:17:02'
$wsPIR.Range('C189').Value = '18:00'
$wsPIR.Range('D189').Value = 'Bathroom'
$wsPIR.Range('E189').Value = 'No Motion'
$wsPIR.Range('F189').Value = 'Inactive'

# ---- Humidity sheet: rows 171-183 ----
$wsHum = $wb.Worksheets.Item("Humidity")
$wsHum.Range('A171').NumberFormat = '@'
$wsHum.Range('A171').Value = '2026-01-28'
$wsHum.Range('B171').Value = '18:16:08'
$wsHum.Range('C171').Value = '18:00'
$wsHum.Range('D171').Value = 'Bathroom'
$wsHum.Range('E171').NumberFormat = '@'
$wsHum.Range('E171').Value = '88.2%'
$wsHum.Range('F171').Value = 'Active'

$wsHum.Range('A172').NumberFormat = '@'
$wsHum.Range('A172').Value = '2026-01-28'
$wsHum.Range('B172').Value = '18:16:11'
$wsHum.Range('C172').Value = '18:00'
$wsHum.Range('D172').Value = 'Bathroom'
$wsHum.Range('E172').NumberFormat = '@'
$wsHum.Range('E172').Value = '88.3%'
$wsHum.Range('F172').Value = 'Active'

$wsHum.Range('A173').NumberFormat = '@'
$wsHum.Range('A173').Value = '2026-01-28'
$wsHum.Range('B173').Value = '18:16:14'
$wsHum.Range('C173').Value = '18:00'
$wsHum.Range('D173').Value = 'Bathroom'
$wsHum.Range('E173').NumberFormat = '@'
$wsHum.Range('E173').Value = '88.2%'
$wsHum.Range('F173').Value = 'Active'

$wsHum.Range('A174').NumberFormat = '@'
$wsHum.Range('A174').Value = '2026-01-28'
$wsHum.Range('B174').Value = '18:16:16'
$wsHum.Range('C174').Value = '18:00'
$wsHum.Range('D174').Value = 'Bathroom'
$wsHum.Range('E174').NumberFormat = '@'
$wsHum.Range('E174').Value = '88.2%'
$wsHum.Range('F174').Value = 'Active'

$wsHum.Range('A175').NumberFormat = '@'
$wsHum.Range('A175').Value = '2026-01-28'
$wsHum.Range('B175').Value = '18:16:20'
$wsHum.Range('C175').Value = '18:00'
$wsHum.Range('D175').Value = 'Bathroom'
$wsHum.Range('E175').NumberFormat = '@'
$wsHum.Range('E175').Value = '87.3%'
$wsHum.Range('F175').Value = 'Active'

$wsHum.Range('A176').NumberFormat = '@'
$wsHum.Range('A176').Value = '2026-01-28'
$wsHum.Range('B176').Value = '18:16:24'
$wsHum.Range('C176').Value = '18:00'
$wsHum.Range('D176').Value = 'Bathroom'
$wsHum.Range('E176').NumberFormat = '@'
$wsHum.Range('E176').Value = '88.3%'
$wsHum.Range('F176').Value = 'Active'

$wsHum.Range('A177').NumberFormat = '@'
$wsHum.Range('A177').Value = '2026-01-28'
$wsHum.Range('B177').Value = '18:16:28'
$wsHum.Range('C177').Value = '18:00'
$wsHum.Range('D177').Value = 'Bathroom'
$wsHum.Range('E177').NumberFormat = '@'
$wsHum.Range('E177').Value = '88.2%'
$wsHum.Range('F177').Value = 'Active'

$wsHum.Range('A178').NumberFormat = '@'
$wsHum.Range('A178').Value = '2026-01-28'
$wsHum.Range('B178').Value = '18:16:32'
$wsHum.Range('C178').Value = '18:00'
$wsHum.Range('D178').Value = 'Bathroom'
$wsHum.Range('E178').NumberFormat = '@'
$wsHum.Range('E178').Value = '87.3%'
$wsHum.Range('F178').Value = 'Active'

$wsHum.Range('A179').NumberFormat = '@'
$wsHum.Range('A179').Value = '2026-01-28'
$wsHum.Range('B179').Value = '18:16:44'
$wsHum.Range('C179').Value = '18:00'
$wsHum.Range('D179').Value = 'Bathroom'
$wsHum.Range('E179').NumberFormat = '@'
$wsHum.Range('E179').Value = '88.2%'
$wsHum.Range('F179').Value = 'Active'

$wsHum.Range('A180').NumberFormat = '@'
$wsHum.Range('A180').Value = '2026-01-28'
$wsHum.Range('B180').Value = '18:16:49'
$wsHum.Range('C180').Value = '18:00'
$wsHum.Range('D180').Value = 'Bathroom'
$wsHum.Range('E180').NumberFormat = '@'
$wsHum.Range('E180').Value = '88.2%'
$wsHum.Range('F180').Value = 'Active'

$wsHum.Range('A181').NumberFormat = '@'
$wsHum.Range('A181').Value = '2026-01-28'
$wsHum.Range('B181').Value = '18:16:53'
$wsHum.Range('C181').Value = '18:00'
$wsHum.Range('D181').Value = 'Bathroom'
$wsHum.Range('E181').NumberFormat = '@'
$wsHum.Range('E181').Value = '87.3%'
$wsHum.Range('F181').Value = 'Active'

$wsHum.Range('A182').NumberFormat = '@'
$wsHum.Range('A182').Value = '2026-01-28'
$wsHum.Range('B182').Value = '18:16:57'
$wsHum.Range('C182').Value = '18:00'
$wsHum.Range('D182').Value = 'Bathroom'
$wsHum.Range('E182').NumberFormat = '@'
$wsHum.Range('E182').Value = '88.2%'
$wsHum.Range('F182').Value = 'Active'

$wsHum.Range('A183').NumberFormat = '@'
$wsHum.Range('A183').Value = '2026-01-28'
$wsHum.Range('B183').Value = '18:17:01'
$wsHum.Range('C183').Value = '18:00'
$wsHum.Range('D183').Value = 'Bathroom'
$wsHum.Range('E183').NumberFormat = '@'
$wsHum.Range('E183').Value = '87.3%'
$wsHum.Range('F183').Value = 'Active'

# ---- Temperature sheet: rows 171-183 ----
$wsTemp = $wb.Worksheets.Item("Temperature")
$wsTemp.Range('A171').NumberFormat = '@'
$wsTemp.Range('A171').Value = '2026-01-28'
$wsTemp.Range('B171').Value = '18:16:09'
$wsTemp.Range('C171').Value = '18:00'
$wsTemp.Range('D171').Value = 'Bathroom'
$wsTemp.Range('E171').Value = '22.9C'
$wsTemp.Range('F171').Value = 'Active'

$wsTemp.Range('A172').NumberFormat = '@'
$wsTemp.Range('A172').Value = '2026-01-28'
$wsTemp.Range('B172').Value = '18:16:12'
$wsTemp.Range('C172').Value = '18:00'
$wsTemp.Range('D172').Value = 'Bathroom'
$wsTemp.Range('E172').Value = '23.0C'
$wsTemp.Range('F172').Value = 'Active'

$wsTemp.Range('A173').NumberFormat = '@'
$wsTemp.Range('A173').Value = '2026-01-28'
$wsTemp.Range('B173').Value = '18:16:15'
$wsTemp.Range('C173').Value = '18:00'
$wsTemp.Range('D173').Value = 'Bathroom'
$wsTemp.Range('E173').Value = '22.9C'
$wsTemp.Range('F173').Value = 'Active'

$wsTemp.Range('A174').NumberFormat = '@'
$wsTemp.Range('A174').Value = '2026-01-28'
$wsTemp.Range('B174').Value = '18:16:17'
$wsTemp.Range('C174').Value = '18:00'
$wsTemp.Range('D174').Value = 'Bathroom'
$wsTemp.Range('E174').Value = '22.9C'
$wsTemp.Range('F174').Value = 'Active'

$wsTemp.Range('A175').NumberFormat = '@'
$wsTemp.Range('A175').Value = '2026-01-28'
$wsTemp.Range('B175').Value = '18:16:21'
$wsTemp.Range('C175').Value = '18:00'
$wsTemp.Range('D175').Value = 'Bathroom'
$wsTemp.Range('E175').Value = '22.9C'
$wsTemp.Range('F175').Value = 'Active'

$wsTemp.Range('A176').NumberFormat = '@'
$wsTemp.Range('A176').Value = '2026-01-28'
$wsTemp.Range('B176').Value = '18:16:25'
$wsTemp.Range('C176').Value = '18:00'
$wsTemp.Range('D176').Value = 'Bathroom'
$wsTemp.Range('E176').Value = '23.0C'
$wsTemp.Range('F176').Value = 'Active'

$wsTemp.Range('A177').NumberFormat = '@'
$wsTemp.Range('A177').Value = '2026-01-28'
$wsTemp.Range('B177').Value = '18:16:29'
$wsTemp.Range('C177').Value = '18:00'
$wsTemp.Range('D177').Value = 'Bathroom'
$wsTemp.Range('E177').Value = '22.9C'
$wsTemp.Range('F177').Value = 'Active'

$wsTemp.Range('A178').NumberFormat = '@'
$wsTemp.Range('A178').Value = '2026-01-28'
$wsTemp.Range('B178').Value = '18:16:33'
$wsTemp.Range('C178').Value = '18:00'
$wsTemp.Range('D178').Value = 'Bathroom'
$wsTemp.Range('E178').Value = '22.9C'
$wsTemp.Range('F178').Value = 'Active'

$wsTemp.Range('A179').NumberFormat = '@'
$wsTemp.Range('A179').Value = '2026-01-28'
$wsTemp.Range('B179').Value = '18:16:45'
$wsTemp.Range('C179').Value = '18:00'
$wsTemp.Range('D179').Value = 'Bathroom'
$wsTemp.Range('E179').Value = '22.9C'
$wsTemp.Range('F179').Value = 'Active'

$wsTemp.Range('A180').NumberFormat = '@'
$wsTemp.Range('A180').Value = '2026-01-28'
$wsTemp.Range('B180').Value = '18:16:49'
$wsTemp.Range('C180').Value = '18:00'
$wsTemp.Range('D180').Value = 'Bathroom'
$wsTemp.Range('E180').Value = '22.9C'
$wsTemp.Range('F180').Value = 'Active'

$wsTemp.Range('A181').NumberFormat = '@'
$wsTemp.Range('A181').Value = '2026-01-28'
$wsTemp.Range('B181').Value = '18:16:54'
$wsTemp.Range('C181').Value = '18:00'
$wsTemp.Range('D181').Value = 'Bathroom'
$wsTemp.Range('E181').Value = '22.9C'
$wsTemp.Range('F181').Value = 'Active'

$wsTemp.Range('A182').NumberFormat = '@'
$wsTemp.Range('A182').Value = '2026-01-28'
$wsTemp.Range('B182').Value = '18:16:58'
$wsTemp.Range('C182').Value = '18:00'
$wsTemp.Range('D182').Value = 'Bathroom'
$wsTemp.Range('E182').Value = '22.9C'
$wsTemp.Range('F182').Value = 'Active'

$wsTemp.Range('A183').NumberFormat = '@'
$wsTemp.Range('A183').Value = '2026-01-28'
$wsTemp.Range('B183').Value = '18:17:02'
$wsTemp.Range('C183').Value = '18:00'
$wsTemp.Range('D183').Value = 'Bathroom'
$wsTemp.Range('E183').Value = '22.9C'
$wsTemp.Range('F183').Value = 'Active'
